$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.979.46'
$ws.Range("E2").Value = '  +2.05%  '
$ws.Range("D3").Value = '2.623.10'
$ws.Range("E3").Value = '  +1.64%  '
$ws.Range("E4").Value = '  +0.05%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '595.93'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.87%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '155.32'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("E7").Value = '  +0.04%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.548'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +1.31%  '
$ws.Range("D9").Value = '2.622.06'
$ws.Range("E9").Value = '  +1.59%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.128'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +11.33%  '
$ws.Range("E11").Value = '  +1.06%  '
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("E13").Value = '  -1.17%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '27.72'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -1.95%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.0000187'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +3.71%  '
$ws.Range("D16").Value = '3.101.91'
$ws.Range("E16").Value = '  +1.30%  '
$ws.Range("D17").Value = '67.751.70'
$ws.Range("E17").Value = '  +1.76%  '
$ws.Range("D18").Value = '2.623.96'
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '368.31'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +3.66%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '11.20'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -0.29%  '
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("E22").Value = '  -0.61%  '
$ws.Range("E23").Value = '  -1.08%  '
$ws.Range("E24").Value = '  +0.21%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '9.87'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -6.77%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '67.45'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.62%  '
$ws.Range("E27").Value = '  +1.52%  '
$ws.Range("D28").Value = '2.730.84'
$ws.Range("E28").Value = '  +0.85%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '577.45'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -5.93%  '
$ws.Range("E30").Value = '  +0.17%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '1.43'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -1.56%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '7.94'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -0.77%  '
$ws.Range("E33").Value = '  +0.59%  '
$ws.Range("E34").Value = '  -1.07%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("E36").Value = '  -3.24%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '4.94'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -2.18%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '159.06'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +2.98%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '19.38'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +0.92%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.371'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.08%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '5.34'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -2.83%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '1.85'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +2.13%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '2.57'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -3.52%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '41.23'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("E45").Value = '  +0.04%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '16.42'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.10%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '156.28'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("E48").Value = '  -6.63%  '
$ws.Range("E49").Value = '  -0.28%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.628'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +2.95%  '
$ws.Range("B51").Value = 'Hedera'
$ws.Range("C51").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.0541'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -3.80%  '
